$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("placesToGo")

# Update the two existing check-in/check-out dates on row 2 (Seattle)
$ws.Range("B2").Value = " 09/10/2022"
$ws.Range("C2").Value = " 11/16/2022"

# Add two new rows of test data. Values are written in this particular
# order so the shared-string table is populated the same way it was
# originally authored.
$ws.Range("A3").Value = "Los Angeles"
$ws.Range("B3").Value = " 12/12/2021"
$ws.Range("A4").Value = "Miami"
$ws.Range("C3").Value = " 12/25/2021"
$ws.Range("C4").Value = " 11/30/2021"
$ws.Range("B4").Value = " 9/28/2021"

# The checkin/checkout columns use the same date-formatted style as the
# existing row. Copy that formatting onto the new cells instead of
# building a brand-new number format.
$ws.Range("B2").Copy()
$ws.Range("B3:C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C5").Select()
